# read_mapping.pptx — "Updating the reproducibility scripts for the new
# version of the manuscript."
#
# 1. Remove the 3rd slide (the BLEND-I / minimap2-Eq comparison figure),
#    which cascades to remove its notes page too.
# 2. Bump the cached "datetimeFigureOut" footer field (11/4/22 -> 11/21/22)
#    everywhere it is stamped: the slide master, every slide layout, and
#    the notes master.

$p = $ppt.ActivePresentation

# --- 1. Delete the third slide (BLEND-I comparison figure) -----------------
$p.Slides.Item(3).Delete()

# --- 2. Refresh the cached date stamp --------------------------------------
function Update-DateStamp($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame -eq $false) { continue }
        $tf = $sh.TextFrame
        if (-not $tf.HasText) { continue }
        $tr = $tf.TextRange
        if ($tr.Text -eq "11/4/22") {
            $tr.Text = "11/21/22"
        }
    }
}

# Slide master
Update-DateStamp $p.SlideMaster.Shapes

# Every slide layout off the master
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DateStamp $layouts.Item($li).Shapes
}

# Notes master
Update-DateStamp $p.NotesMaster.Shapes
